# Update "想去人数" (F column) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 191
    4  = 12609
    5  = 1282
    6  = 151
    10 = 207
    16 = 384
    17 = 4318
    18 = 99
    19 = 23
    20 = 947
    21 = 24
    23 = 79
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
